# Updates the cryptos price/volume table to the latest scraped values.
# Column D ("Price") holds text-formatted numbers (not real numbers) in the
# original workbook, so any value that Excel would otherwise auto-convert to
# a number needs to be written with a text number format to keep it a string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

function Set-Value($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "35.359.88"
Set-Value      "E2" "  +1.80%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.882.59"
Set-Value      "E3" "  +0.45%  "

# Row 4 - TetherUSD
Set-Value "E4" "  -0.10%  "

# Row 5 - was XRP, now BNB
Set-Value     "B5" "BNB"
Set-Value     "C5" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue "D5" "245.51"
Set-Value     "E5" "  -0.82%  "

# Row 6 - was BNB, now XRP
Set-Value     "B6" "XRP"
Set-Value     "C6" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue "D6" "0.690"
Set-Value     "E6" "  +0.21%  "

# Row 7 - USDC
Set-Value "E7" "  -0.03%  "

# Row 8 - Solana
Set-TextValue "D8" "43.18"
Set-Value      "E8" "  +2.93%  "

# Row 9 - Cardano
Set-Value "E9" "  +2.05%  "

# Row 10 - OKB
Set-TextValue "D10" "54.88"
Set-Value      "E10" "  +7.19%  "

# Row 11 - Dogecoin
Set-Value "E11" "  +1.01%  "

# Row 12 - TRON
Set-Value "E12" "  +0.96%  "

# Row 13 - Chainlink
Set-TextValue "D13" "13.69"
Set-Value      "E13" "  +6.99%  "

# Row 14 - WrappedliquidstakedEther2.0 (only price changes)
Set-TextValue "D14" "2.157.11"

# Row 15 - Polygon
Set-TextValue "D15" "0.768"
Set-Value      "E15" "  +7.70%  "

# Row 16 - Polkadot
Set-TextValue "D16" "5.01"
Set-Value      "E16" "  +2.20%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "1.888.80"
Set-Value      "E17" "  +0.71%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "35.348.18"
Set-Value      "E18" "  +1.75%  "

# Row 19 - Litecoin
Set-TextValue "D19" "73.28"
Set-Value      "E19" "  +0.72%  "

# Row 20 - ShibaInu
Set-TextValue "D20" "0.0₃0823"
Set-Value      "E20" "  +0.61%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "243.76"
Set-Value      "E21" "  -0.14%  "

# Row 22 - Avalanche
Set-Value "E22" "  +0.87%  "

# Row 23 - Uniswap
Set-TextValue "D23" "5.12"
Set-Value      "E23" "  +4.45%  "

# Row 24 - Toncoin
Set-TextValue "D24" "2.63"
Set-Value      "E24" "  +8.60%  "

# Row 25 - Dai
Set-Value "E25" "  +0.11%  "

# Row 26 - PancakeSwap
Set-TextValue "D26" "2.15"
Set-Value      "E26" "  -4.12%  "

# Row 27 - Monero
Set-TextValue "D27" "165.43"
Set-Value      "E27" "  +0.30%  "

# Row 28 - Cosmos
Set-TextValue "D28" "8.58"
Set-Value      "E28" "  +2.57%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "18.23"
Set-Value      "E29" "  +0.37%  "

# Row 30 - Stellar
Set-Value "E30" "  -0.11%  "

# Row 31 - Hedera
Set-Value "E31" "  +3.27%  "

# Row 32 - Filecoin
Set-TextValue "D32" "4.28"
Set-Value      "E32" "  -0.38%  "

# Row 33 - WEMIXToken
Set-TextValue "D33" "1.88"
Set-Value      "E33" "  +15.01%  "

# Row 34 - InternetComputer(DFINITY)
Set-Value "E34" "  -0.22%  "

# Row 35 - BinanceUSD
Set-Value "E35" "  -0.12%  "

# Row 36 - TrustWalletToken
Set-TextValue "D36" "1.46"
Set-Value      "E36" "  -12.23%  "

# Row 37 - ImmutableX
Set-TextValue "D37" "0.847"
Set-Value      "E37" "  +2.19%  "

# Row 38 - LidoDAOToken
Set-TextValue "D38" "1.92"
Set-Value      "E38" "  -2.77%  "

# Row 39 - Kaspa
Set-TextValue "D39" "0.0718"
Set-Value      "E39" "  +8.65%  "

# Row 40 - VeChain
Set-TextValue "D40" "0.0219"
Set-Value      "E40" "  +4.26%  "

# Row 41 - was Aave, now InjectiveProtocol
Set-Value     "B41" "InjectiveProtocol"
Set-Value     "C41" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D41" "17.05"
Set-Value     "E41" "  +1.11%  "

# Row 42 - was InjectiveProtocol, now Aave
Set-Value     "B42" "Aave"
Set-Value     "C42" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D42" "97.23"
Set-Value     "E42" "  -0.66%  "

# Row 43 - ARBITRUM
Set-Value "E43" "  -1.03%  "

# Row 44 - Gas
Set-TextValue "D44" "13.66"
Set-Value      "E44" "  +11.76%  "

# Row 45 - Maker
Set-TextValue "D45" "1.312.55"
Set-Value      "E45" "  +2.25%  "

# Row 46 - RenderToken
Set-TextValue "D46" "2.38"
Set-Value      "E46" "  +1.76%  "

# Row 47 - Cronos
Set-TextValue "D47" "0.0811"
Set-Value      "E47" "  +3.92%  "

# Row 48 - HuobiToken
Set-Value "E48" "  +0.02%  "

# Row 49 - MXToken
Set-Value "E49" "  +0.38%  "

# Row 50 - FraxShare
Set-TextValue "D50" "6.26"
Set-Value      "E50" "  -2.97%  "

# Row 51 - RocketPoolETH
Set-TextValue "D51" "2.056.46"
Set-Value      "E51" "  -0.06%  "
